# Generate Report for Handoff
# Updates the localization-status report: marks items as "Ready for handoff"
# and refreshes the handoff timestamps, then narrows the now-shorter
# status/date columns on each sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update status text from "Handed back: in sync with en-US" to "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Refresh handoff timestamps (stored as plain text) ---
$overview.Range("G2").Value = "2016-09-07 03:15:42"
$dede.Range("H2").Value     = "2016-09-07 03:15:42"
$zhcn.Range("H2").Value     = "2016-09-07 03:15:36"

# --- Narrow the status/date columns to match the shorter content ---
# (target stored width ~17.2159881591797; ColumnWidth snaps to the nearest
# pixel-quantized grid value, so feed it the input that lands closest.)
$overview.Columns.Item(5).ColumnWidth = 16.3333333333333
$overview.Columns.Item(6).ColumnWidth = 16.3333333333333
$zhcn.Columns.Item(3).ColumnWidth     = 16.3333333333333
$dede.Columns.Item(3).ColumnWidth     = 16.3333333333333
